# Fix the header row "syntax error": rename UPPER_SNAKE_CASE column headers
# to camelCase field names.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "txnId"
$ws.Range("B1").Value = "acctId"
$ws.Range("C1").Value = "txnSeq"
$ws.Range("D1").Value = "txnType"
$ws.Range("E1").Value = "tradeDate"
$ws.Range("F1").Value = "fundId"
$ws.Range("G1").Value = "valnDate"
$ws.Range("H1").Value = "unit"
$ws.Range("I1").Value = "processDate"
$ws.Range("J1").Value = "unitCost"

# Move the active selection to J1, matching the saved selection state.
$ws.Range("J1").Select()
